$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tables")
$ws.Activate()

# Sales order: add "payment_type" field before "progress"
$ws.Range("B14").Value = "docno,date,user,customer,fycode,fncode,payment_type,progress,_ref"

# Sales order items: add "tax,discount,total" fields before "_ref"
$ws.Range("B15").Value = "so,product,rate,quantity,tax,discount,total,_ref"

# Update selection state to match post-edit state
$ws.Range("B16").Select()
